{"js": "// Add a visible \"Date of Submission: {date}\" line, keeping the original\n// (white / invisible-text) date stamp run, whose own text is trimmed\n// from \"2023-03-30\" to \"3-30\". The \"_GoBack\" bookmark that used to sit\n// right before the old single run now sits right before that trimmed\n// white run (i.e. after the two newly-added runs).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that holds the submission date line.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Date of Submission:\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the 'Date of Submission:' paragraph\");\n}\n\n// 1) \"Date of Submission: \" (trailing space) -> \"Date of Submission:\"\nconst labelHits = target.search(\"Date of Submission: \", { matchCase: true });\nlabelHits.load(\"items\");\nawait context.sync();\nif (labelHits.items.length === 0) {\n  throw new Error(\"Could not find the 'Date of Submission: ' run\");\n}\nconst labelRun = labelHits.items[0];\nlabelRun.insertText(\"Date of Submission:\", \"Replace\");\nawait context.sync();\n\n// 2) Insert the new \" {date}\" placeholder text right after the label,\n//    re-using the label run's own (non-colored) formatting.\nconst freshLabelHits = target.search(\"Date of Submission:\", { matchCase: true });\nfreshLabelHits.load(\"items\");\nawait context.sync();\nconst freshLabelRun = freshLabelHits.items[0];\nfreshLabelRun.insertText(\" {date}\", \"After\");\nawait context.sync();\n\n// 3) Move the \"_GoBack\" bookmark so it now sits right before the old\n//    (white / hidden) date-stamp run instead of right before the label.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst oldDateHits = target.search(\"2023-03-30\", { matchCase: true });\noldDateHits.load(\"items\");\nawait context.sync();\nif (oldDateHits.items.length === 0) {\n  throw new Error(\"Could not find the '2023-03-30' run\");\n}\nconst oldDateStart = oldDateHits.items[0].getRange(\"Start\");\noldDateStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) \"2023-03-30\" -> \"3-30\" (bugfix mentioned in the commit message).\nconst oldDateHits2 = target.search(\"2023-03-30\", { matchCase: true });\noldDateHits2.load(\"items\");\nawait context.sync();\nconst oldDateRun2 = oldDateHits2.items[0];\noldDateRun2.insertText(\"3-30\", \"Replace\");\nawait context.sync();\n", "ps1": "# Add a visible \"Date of Submission: {date}\" line, keeping the original\n# (white / invisible-text) date stamp run, whose own text is trimmed\n# from \"2023-03-30\" to \"3-30\". The \"_GoBack\" bookmark that used to sit\n# right before the old single run now sits right before that trimmed\n# white run (i.e. after the two newly-added runs).\n\n$d = $word.ActiveDocument\n\n# 1) \"Date of Submission: \" (trailing space) -> \"Date of Submission:\"\n$labelRange = $d.Content\n$labelRange.Find.Execute(\"Date of Submission: \", $true) | Out-Null\n$labelRange.Text = \"Date of Submission:\"\n\n# 2) Insert the new \" {date}\" placeholder text right after the label,\n#    re-using the label run's own (non-colored) formatting.\n$labelRange2 = $d.Content\n$labelRange2.Find.Execute(\"Date of Submission:\", $true) | Out-Null\n$labelRange2.InsertAfter(\" {date}\")\n\n# 3) Move the \"_GoBack\" bookmark so it now sits right before the old\n#    (white / hidden) date-stamp run instead of right before the label.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Delete()\n\n$dateRange = $d.Content\n$dateRange.Find.Execute(\"2023-03-30\", $true) | Out-Null\n$dateRange.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $dateRange)\n\n# 4) \"2023-03-30\" -> \"3-30\" (bugfix mentioned in the commit message).\n$oldDateRange = $d.Content\n$oldDateRange.Find.Execute(\"2023-03-30\", $true) | Out-Null\n$oldDateRange.Text = \"3-30\"\n"}
